$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old 2-column table entirely (values + formatting) ---
$ws.Range("A1:B10").Clear()

# --- New column widths: column A wider (label column), column B narrower ---
# (ColumnWidth undergoes pixel-rounding internally; these inputs land exactly
#  on stored widths of 38 and 10.5 respectively.)
$ws.Columns.Item(1).ColumnWidth = 37.166666666666664
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666

# --- Header row (row 2): Weights | 31 - 35 | 36 - 40 | ... | 71 - 75 ---
$ws.Range("A2").Value = "Weights"
$ws.Range("A2").Font.Bold = $true

$ws.Range("B2").Value = "31 - 35"
$ws.Range("C2").Value = "36 - 40"
$ws.Range("D2").Value = "41 - 45"
$ws.Range("E2").Value = "46 - 50"
$ws.Range("F2").Value = "51 - 55"
$ws.Range("G2").Value = "56 - 60"
$ws.Range("H2").Value = "61 - 65"
$ws.Range("I2").Value = "66 - 70"
$ws.Range("J2").Value = "71 - 75"

# --- Data row (row 3): label + counts ---
$ws.Range("A3").Value = "Number of  Students (frequency)"
$ws.Range("A3").Font.Bold = $true

$ws.Range("B3").Value = 9
$ws.Range("B3").HorizontalAlignment = -4108

$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1

# --- Empty placeholder rows below the table (rows 4-11, columns A:B) ---
$ws.Range("A4:B11").HorizontalAlignment = -4108

# --- Selection / active cell ---
$ws.Range("A4").Select()
